$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Narrow column D slightly (stored width 13 -> 12).
# Note: Excel's ColumnWidth property (chars) differs from the stored OOXML
# width by a constant offset (~0.8333333333333333) on this engine, so we
# compensate to land on an exact stored width of 12.
$ws.Columns.Item(4).ColumnWidth = 12 - 0.8333333333333333

# Update VENTA (D) and POR CUMPLIR (E) for the "OTROS" row (row 2)
$ws.Range("D2").Value = 3076.4
$ws.Range("E2").Value = -3076.4

# Update TOTAL row (row 4): VENTA, POR CUMPLIR, CUMPLIMIENTO
$ws.Range("D4").Value = 3324.88
$ws.Range("E4").Value = 14175.12
$ws.Range("F4").Value = 0.1899931428571429
